$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "7. Do we need to number abbreviations as reference numbers"
$ws.Range("A8").Select()
